# Apply crypto price/volume updates per commit "Updated cryptos list" diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.660.18"
$ws.Range("E2").Value = "  -0.66%  "

$ws.Range("D3").Value = "'2.543.67"
$ws.Range("E3").Value = "  +0.25%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "'311.61"
$ws.Range("E5").Value = "  -1.75%  "

$ws.Range("D6").Value = "'101.73"
$ws.Range("E6").Value = "  +5.51%  "

$ws.Range("D7").Value = "'0.570"
$ws.Range("E7").Value = "  -1.07%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("E9").Value = "  -1.37%  "

$ws.Range("D10").Value = "'36.06"
$ws.Range("E10").Value = "  +1.27%  "

$ws.Range("D11").Value = "'0.0805"
$ws.Range("E11").Value = "  -1.43%  "

$ws.Range("E12").Value = "  -1.42%  "

$ws.Range("E13").Value = "  -0.55%  "

$ws.Range("D14").Value = "'2.938.33"
$ws.Range("E14").Value = "  +0.27%  "

$ws.Range("D15").Value = "'15.96"
$ws.Range("E15").Value = "  +6.13%  "

$ws.Range("D16").Value = "'2.546.19"
$ws.Range("E16").Value = "  -1.71%  "

$ws.Range("E17").Value = "  -1.82%  "

$ws.Range("D18").Value = "'42.673.33"
$ws.Range("E18").Value = "  -0.82%  "

$ws.Range("D19").Value = "'6.83"
$ws.Range("E19").Value = "  +0.46%  "

$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "'0.0₃0954"
$ws.Range("E20").Value = "  -1.21%  "

$ws.Range("B21").Value = "InternetComputer(DFINITY)"
$ws.Range("C21").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D21").Value = "'12.31"
$ws.Range("E21").Value = "  -1.91%  "

$ws.Range("D22").Value = "'69.14"
$ws.Range("E22").Value = "  -0.89%  "

$ws.Range("D23").Value = "'244.43"
$ws.Range("E23").Value = "  -3.55%  "

$ws.Range("E24").Value = "  -1.00%  "

$ws.Range("E25").Value = "  +0.06%  "

$ws.Range("D26").Value = "'26.50"
$ws.Range("E26").Value = "  -1.55%  "

$ws.Range("E27").Value = "  +0.01%  "

$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").Value = "'2.37"
$ws.Range("E28").Value = "  -1.15%  "

$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").Value = "'40.61"
$ws.Range("E29").Value = "  -0.26%  "

$ws.Range("E30").Value = "  -2.51%  "

$ws.Range("D31").Value = "'157.90"
$ws.Range("E31").Value = "  +1.22%  "

$ws.Range("E32").Value = "  -3.20%  "

$ws.Range("D33").Value = "'2.76"
$ws.Range("E33").Value = "  +12.27%  "

$ws.Range("E34").Value = "  +0.77%  "

$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "'2.06"
$ws.Range("E35").Value = "  -2.32%  "

$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").Value = "'2.63"
$ws.Range("E36").Value = "  -3.25%  "

$ws.Range("E37").Value = "  -3.56%  "

$ws.Range("D38").Value = "'18.38"
$ws.Range("E38").Value = "  -5.12%  "

$ws.Range("D39").Value = "'0.111"
$ws.Range("E39").Value = "  -1.87%  "

$ws.Range("E40").Value = "  -0.63%  "

$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "'4.19"
$ws.Range("E41").Value = "  +9.72%  "

$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").Value = "'22.31"
$ws.Range("E42").Value = "  +2.34%  "

$ws.Range("E43").Value = "  +0.15%  "

$ws.Range("D44").Value = "'3.33"
$ws.Range("E44").Value = "  +1.49%  "

$ws.Range("E45").Value = "  -1.76%  "

$ws.Range("D46").Value = "'1.975.06"
$ws.Range("E46").Value = "  -1.26%  "

$ws.Range("D47").Value = "'8.90"
$ws.Range("E47").Value = "  -1.83%  "

$ws.Range("D48").Value = "'2.791.67"
$ws.Range("E48").Value = "  +0.21%  "

$ws.Range("D49").Value = "'81.30"
$ws.Range("E49").Value = "  -3.88%  "

$ws.Range("E50").Value = "  +0.44%  "

$ws.Range("D51").Value = "'73.28"
$ws.Range("E51").Value = "  -2.09%  "
